{"js": "// Update the answer cells in the two-digit division table.\n// Each entry is a unique \"before\" text that is replaced in-place with\n// the corresponding \"after\" text, mirroring the canonical OOXML diff.\nconst replacements = [\n  [\"79\u00f78=9, 7\", \"42\u00f75=8, 2\"],\n  [\"73\u00f73=24, 1\", \"48\u00f72=24, 0\"],\n  [\"21\u00f77=3, 0\", \"74\u00f73=24, 2\"],\n  [\"33\u00f72=16, 1\", \"52\u00f73=17, 1\"],\n  [\"26\u00f76=4, 2\", \"20\u00f78=2, 4\"],\n  [\"45\u00f76=7, 3\", \"59\u00f76=9, 5\"],\n  [\"37\u00f78=4, 5\", \"25\u00f77=3, 4\"],\n  [\"67\u00f73=22, 1\", \"12\u00f73=4, 0\"],\n  [\"41\u00f78=5, 1\", \"37\u00f73=12, 1\"],\n  [\"36\u00f75=7, 1\", \"17\u00f79=1, 8\"],\n  [\"99\u00f78=12, 3\", \"78\u00f75=15, 3\"],\n  [\"61\u00f72=30, 1\", \"99\u00f73=33, 0\"],\n  [\"78\u00f73=26, 0\", \"23\u00f74=5, 3\"],\n  [\"23\u00f79=2, 5\", \"94\u00f79=10, 4\"],\n  [\"88\u00f72=44, 0\", \"90\u00f78=11, 2\"],\n  [\"65\u00f78=8, 1\", \"15\u00f76=2, 3\"],\n  [\"68\u00f78=8, 4\", \"19\u00f78=2, 3\"],\n  [\"73\u00f77=10, 3\", \"29\u00f76=4, 5\"],\n  [\"96\u00f72=48, 0\", \"55\u00f78=6, 7\"],\n  [\"80\u00f77=11, 3\", \"86\u00f73=28, 2\"],\n  [\"99\u00f79=11, 0\", \"64\u00f77=9, 1\"],\n  [\"51\u00f79=5, 6\", \"25\u00f79=2, 7\"],\n  [\"20\u00f76=3, 2\", \"56\u00f79=6, 2\"],\n  [\"78\u00f78=9, 6\", \"48\u00f73=16, 0\"],\n  [\"70\u00f72=35, 0\", \"56\u00f79=6, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the answer cells in the two-digit division table.\n# Each entry is a unique \"before\" text that is replaced in-place with\n# the corresponding \"after\" text, mirroring the canonical OOXML diff.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"79\u00f78=9, 7\", \"42\u00f75=8, 2\"),\n    @(\"73\u00f73=24, 1\", \"48\u00f72=24, 0\"),\n    @(\"21\u00f77=3, 0\", \"74\u00f73=24, 2\"),\n    @(\"33\u00f72=16, 1\", \"52\u00f73=17, 1\"),\n    @(\"26\u00f76=4, 2\", \"20\u00f78=2, 4\"),\n    @(\"45\u00f76=7, 3\", \"59\u00f76=9, 5\"),\n    @(\"37\u00f78=4, 5\", \"25\u00f77=3, 4\"),\n    @(\"67\u00f73=22, 1\", \"12\u00f73=4, 0\"),\n    @(\"41\u00f78=5, 1\", \"37\u00f73=12, 1\"),\n    @(\"36\u00f75=7, 1\", \"17\u00f79=1, 8\"),\n    @(\"99\u00f78=12, 3\", \"78\u00f75=15, 3\"),\n    @(\"61\u00f72=30, 1\", \"99\u00f73=33, 0\"),\n    @(\"78\u00f73=26, 0\", \"23\u00f74=5, 3\"),\n    @(\"23\u00f79=2, 5\", \"94\u00f79=10, 4\"),\n    @(\"88\u00f72=44, 0\", \"90\u00f78=11, 2\"),\n    @(\"65\u00f78=8, 1\", \"15\u00f76=2, 3\"),\n    @(\"68\u00f78=8, 4\", \"19\u00f78=2, 3\"),\n    @(\"73\u00f77=10, 3\", \"29\u00f76=4, 5\"),\n    @(\"96\u00f72=48, 0\", \"55\u00f78=6, 7\"),\n    @(\"80\u00f77=11, 3\", \"86\u00f73=28, 2\"),\n    @(\"99\u00f79=11, 0\", \"64\u00f77=9, 1\"),\n    @(\"51\u00f79=5, 6\", \"25\u00f79=2, 7\"),\n    @(\"20\u00f76=3, 2\", \"56\u00f79=6, 2\"),\n    @(\"78\u00f78=9, 6\", \"48\u00f73=16, 0\"),\n    @(\"70\u00f72=35, 0\", \"56\u00f79=6, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n}\n"}
